$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '328.84'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '0.50%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '44.22'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '1.18%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.578'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '1.68%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.08071'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '0.11%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.971'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '5.11%'
$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '4.327'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '0.84%'
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9528'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '1.76%'
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1160'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-2.35%'
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1854'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-1.84%'
$ws.Range("B12").Value = 'MCDex'
$ws.Range("C12").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '11.82'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '37.77%'
$ws.Range("B13").Value = 'MandalaExchangeToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09787'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '2.88%'
$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.04682'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '13.61%'
$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.1067'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.07%'
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001284'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '0.20%'
$ws.Range("B17").Value = 'CoinExToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.04234'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-2.46%'
$ws.Range("B18").Value = 'TigerCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.005946'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-0.71%'
$ws.Range("B19").Value = 'HotbitToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.004319'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-0.30%'
$ws.Range("B20").Value = 'LEO'
$ws.Range("C20").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.372'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '-5.71%'
$ws.Range("B21").Value = 'BitpandaEcosystemToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.3474'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-0.31%'
$ws.Range("B22").Value = 'ProBitToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.1409'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '3.20%'
$ws.Range("B23").Value = 'ZBToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.2507'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-3.30%'
$ws.Range("B24").Value = 'BitKan'
$ws.Range("C24").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001252'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '1.36%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0001191'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-3.42%'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-0.51%'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02636'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '0.01%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05540'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '2.67%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.007554'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-1.51%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1408'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '1.46%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.008087'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-27.07%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002017'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-5.20%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008900'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-7.90%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00007131'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '3.76%'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-0.17%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.002301'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '1.12%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.003526'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-1.19%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002102'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-0.17%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0002002'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.17%'
